$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B228').Value = 'IgE ESPECIFICO (E82) - Epitélios - PELO DE COELHO'
$ws.Range('B252').Value = 'IgE ESPECIFICO (K70) - Ocupacionais - GRAO DE CAFE'
$ws.Range('B290').Value = 'PARASITOLOGICO 1ª AMOSTRA'
$ws.Range('B291').Value = 'PARASITOLOGICO 2ª AMOSTRA'
$ws.Range('B292').Value = 'PARASITOLOGICO 3ª AMOSTRA'
$ws.Range('B293').Value = 'PARASITOLOGICO DE FEZES'
$ws.Range('B419').Value = 'RM PESCOÇO'
$ws.Range('B456').Value = 'ROTAVIRUS - Pesquisa'
$ws.Range('B457').Value = 'RUBEOLA - Anticorpos IgG'
$ws.Range('B458').Value = 'RUBEOLA - Anticorpos IgM'
$ws.Range('B532').Value = 'RX PESCOÇO'
$ws.Range('B535').Value = 'RX PE D'
$ws.Range('B536').Value = 'RX PE E'
$ws.Range('B551').Value = 'SELENIO SERICO'
$ws.Range('B555').Value = 'SUMARIO DE URINA'
$ws.Range('B556').Value = 'SIFILIS - VDRL'
$ws.Range('B557').Value = 'SODIO'
$ws.Range('B558').Value = 'SODIO URINARIO - 24h'
$ws.Range('B559').Value = 'SODIO URINARIO - AMOSTRA ISOLADA'
$ws.Range('B574').Value = 'TC ARTICULAÇAO EXTERNO CLAVICULAR D'
$ws.Range('B575').Value = 'TC ARTICULAÇAO EXTERNO CLAVICULAR E'
$ws.Range('B576').Value = 'TC ARTICULAÇOES MEMBROS INFERIORES'
$ws.Range('B604').Value = 'TC HIPOFISE C/CONTRASTE'
$ws.Range('B605').Value = 'TC HIPOFISE S/CONTRASTE'
$ws.Range('B664').Value = 'TESTE DE TOLERANCIA A GLICOSE'
$ws.Range('B665').Value = 'TESTE ERGOMETRICO'
$ws.Range('B681').Value = 'TRIGLICERIDEOS'
$ws.Range('B686').Value = 'UREIA URINARIA - 24h'
$ws.Range('B693').Value = 'USG ENDOVAGINAL + PELVICA'
$ws.Range('B694').Value = 'USG ENDOVAGINAL + PELVICA DOPPLER'
$ws.Range('B700').Value = 'USG PELVICA'
$ws.Range('B701').Value = 'USG PELVICA DOPPLER'
$ws.Range('B703').Value = 'USG VIAS URINARIAS'
$ws.Range('B704').Value = 'USG VIAS URINARIAS DOPPLER'
$ws.Range('B707').Value = 'VASOPRESSINA - antidiuretico - AD'
$ws.Range('B715').Value = 'ZINCO SERICO'
$ws.Range('B716').Value = 'ACIDO LATICO'
$ws.Range('B717').Value = 'ACIDO VANIL MANDELICO'
$ws.Range('B718').Value = 'ACIDO URICO'
$ws.Range('B719').Value = 'INDICE DE SATURAÇÃO DA TRANSFERRINA'
$ws.Range('B720').Value = 'INDICE TIROXINA LIVRE'

$win = $excel.ActiveWindow
$win.ScrollRow = 168
$win.ScrollColumn = 1

